$wb = $excel.ActiveWorkbook

# Sheets affected: zh-cn (row 7) and de-de (row 7).
# Both receive the same kind of update: the handback for
# 4ffced0e-da37-42ac-8f6d-07ff54a8e8fc is complete but stale, so the
# "Latest Target/Handback File", "Latest Handback DateTime" and
# "Error Detail" columns for that row get populated, plus a hyperlink
# on the newly filled "Latest Target File" cell (column I).

$sheets = @(
    @{ Name = "zh-cn"; Xlf = "4ffced0e-da37-42ac-8f6d-07ff54a8e8fc.d1d6951112ff798d2a8c0d5fd3f63cad69845773.zh-cn.xlf"; DateTime = "2016-08-12 09:06:55" },
    @{ Name = "de-de"; Xlf = "4ffced0e-da37-42ac-8f6d-07ff54a8e8fc.d1d6951112ff798d2a8c0d5fd3f63cad69845773.de-de.xlf"; DateTime = "2016-08-12 09:07:08" }
)

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/7271ef1490b55a09418739690a2a168bc389f6be/e2e/4ffced0e-da37-42ac-8f6d-07ff54a8e8fc.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/55474f1393e873da5bba8d8a0867a1c9214a9302/e2e/4ffced0e-da37-42ac-8f6d-07ff54a8e8fc.md."

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # I7: Latest Target File -> becomes a hyperlink to the source .md file
    $ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/55474f1393e873da5bba8d8a0867a1c9214a9302/e2e/4ffced0e-da37-42ac-8f6d-07ff54a8e8fc.md", $null, $null, "4ffced0e-da37-42ac-8f6d-07ff54a8e8fc.md")

    # J7: Latest Handback File
    $ws.Range("J7").Value = $info.Xlf

    # K7: Latest Handback DateTime
    $ws.Range("K7").Value = $info.DateTime

    # P7: Error Detail
    $ws.Range("P7").Value = $errorDetail
}
